# Sprint Backlog: mark the remaining open tasks as Done.
#
# The tracker's Status column (G5:H15, validated against
# "Not yet started,In Progress,Done") still had a few tasks flagged
# "In Progress" (orange fill) or "Not yet started" (red fill). The sprint
# wrapped up, so every task is now complete: set those cells' text to
# "Done" and recolor them to match the existing "Done" rows' green fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matches the fill already used on the other "Done" status cells (RGB 00B050,
# i.e. OLE_COLOR 0x0050B000 = B*65536 + G*256 + R).
$doneGreen = 5287936

$statusRanges = @("G10:H11", "G13:H13", "G14:H15")
foreach ($rangeAddr in $statusRanges) {
    $rng = $ws.Range($rangeAddr)
    $rng.Value2 = "Done"
    $rng.Interior.Color = $doneGreen
}

# Reflect the final review pass: cursor left resting on the disclaimer note.
$ws.Range("C19:N19").Select() | Out-Null
